# feat: add 2022-Q1 data
#
# The former "总计" (Total) sheet (5th sheet) becomes the new "2022-Q1"
# quarterly holdings sheet, and a brand-new "总计" sheet is appended at the
# end containing the historical summary table plus a new row for 2022-Q1.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Turn the old "总计" sheet into the new "2022-Q1" quarterly sheet.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item(5)
$q1.Name = "2022-Q1"

# Clear out the old summary-table content (B1:D5) before laying out the
# new quarterly holdings table. ClearContents (not Clear) so the header /
# index-column styling (bold + border on row 1 and column A) survives.
$q1.Range("A1:D5").ClearContents()

# The header row needs to grow from B1:D1 out to B1:H1 and the index
# column from A2:A5 down to A2:A7 - copy the existing (already-styled)
# header/index cells' formatting onto the newly-needed cells first so the
# whole row/column keeps the same bold+border+centered look.
$q1.Cells.Item(1, 2).Copy() | Out-Null
$q1.Range("E1:H1").PasteSpecial(-4122) | Out-Null
$q1.Cells.Item(2, 1).Copy() | Out-Null
$q1.Range("A6:A7").PasteSpecial(-4122) | Out-Null

# Header row
$q1.Cells.Item(1, 2).Value = "基金代码"
$q1.Cells.Item(1, 3).Value = "基金名称"
$q1.Cells.Item(1, 4).Value = "基金规模"
$q1.Cells.Item(1, 5).Value = "股票总仓位"
$q1.Cells.Item(1, 6).Value = "仓位占比"
$q1.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q1.Cells.Item(1, 8).Value = "仓位排名"

# Fund rows: index, code, name, scale, stock position, position share,
# held market value (CNY 100M), position rank.
$q1Rows = @(
    @(0, "012071", "中加喜利回报一年持有期混合A", "5.20", "46.57", "3.09", "0.1607", 3),
    @(1, "005775", "中加转型动力灵活配置混合A",   "3.41", "66.34", "3.92", "0.1337", 4),
    @(2, "009242", "中加核心智造混合A",           "2.05", "65.71", "5.06", "0.1037", 2),
    @(3, "012072", "中加喜利回报一年持有期混合C", "2.74", "46.57", "3.09", "0.0847", 3),
    @(4, "005776", "中加转型动力灵活配置混合C",   "1.92", "66.34", "3.92", "0.0753", 4),
    @(5, "009243", "中加核心智造混合C",           "0.10", "65.71", "5.06", "0.0051", 2)
)

# Force text so numeric-looking strings (fund codes, scale, percentages)
# keep their exact printed form (leading/trailing zeros) instead of being
# coerced to numbers.
$q1.Range("B2:G7").NumberFormat = "@"

$r = 2
foreach ($row in $q1Rows) {
    $q1.Cells.Item($r, 1).Value = $row[0]
    $q1.Cells.Item($r, 2).Value = $row[1]
    $q1.Cells.Item($r, 3).Value = $row[2]
    $q1.Cells.Item($r, 4).Value = $row[3]
    $q1.Cells.Item($r, 5).Value = $row[4]
    $q1.Cells.Item($r, 6).Value = $row[5]
    $q1.Cells.Item($r, 7).Value = $row[6]
    $q1.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2) Append a brand-new "总计" sheet after "2022-Q1" with the historical
#    summary table (date / holding count / holding value), now including
#    the 2022-Q1 row at the top.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

# Match the look & feel (page margins / sheetPr) of the other data sheets.
$total.PageSetup.LeftMargin = 0.75 * 72
$total.PageSetup.RightMargin = 0.75 * 72
$total.PageSetup.TopMargin = 1 * 72
$total.PageSetup.BottomMargin = 1 * 72
$total.PageSetup.HeaderMargin = 0.5 * 72
$total.PageSetup.FooterMargin = 0.5 * 72

# Pull the bold+border+centered header/index-column formatting over from
# the 2022-Q1 sheet (same look as every other sheet's header row / column
# A) since a brand-new sheet starts out with no styling at all.
$q1.Cells.Item(1, 2).Copy() | Out-Null
$total.Range("B1:D1").PasteSpecial(-4122) | Out-Null
$q1.Cells.Item(2, 1).Copy() | Out-Null
$total.Range("A2:A6").PasteSpecial(-4122) | Out-Null

$total.Cells.Item(1, 2).Value = "日期"
$total.Cells.Item(1, 3).Value = "持有数量(只)"
$total.Cells.Item(1, 4).Value = "持有市值(亿元)"

$totalRows = @(
    @(0, "2022-Q1", 6, 0.5600000000000001),
    @(1, "2021-Q4", 13, 0.92),
    @(2, "2021-Q3", 4, 0.33),
    @(3, "2021-Q2", 9, 0.88),
    @(4, "2021-Q1", 8, 0.24)
)

$r = 2
foreach ($row in $totalRows) {
    $total.Cells.Item($r, 1).Value = $row[0]
    $total.Cells.Item($r, 2).Value = $row[1]
    $total.Cells.Item($r, 3).Value = $row[2]
    $total.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 3) Restore the original active sheet/selection (2021-Q1) so adding the
#    new sheet doesn't shift which tab is selected.
# ---------------------------------------------------------------------
$wb.Worksheets.Item(1).Activate()
